$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "questions = [`n    {`n        `"title`": `"You have received the following schedule from the warehouse planning team. A picker can pick five orders in a two-hour slot. How many orders have been picked by lunchtime?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"30`",`n            `"50`",`n            `"60`",`n            `"100`"`n        ],`n        `"score`": `"50`"`n    },`n    {`n        `"title`": `"You requested opening stock numbers from the warehouse IT department, but the data entry operator forgot to input the opening stock of \u201cTables.\u201d The table below summarizes the net receipt in a two-hour window, where net receipt = total receipt - total dispatch.  What was the opening stock of \u201cTables\u201d?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"0`",`n            `"12`",`n            `"21`",`n            `"30`"`n        ],`n        `"score`": `"12`"`n    },`n    {`n        `"title`": `"It is May 31. The warehouse where you work uses the following table to track equipment maintenance. You are filling in the maintenance due dates that are missing. Any equipment requiring emergency maintenance will undergo maintenance ahead of its due date and within the next 30 days. How many pieces of equipment will undergo maintenance in June?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"2`",`n            `"3`",`n            `"5`",`n            `"8`"`n        ],`n        `"score`": `"5`"`n    },`n    {`n        `"title`": `"Your warehouse assesses the quality score of each picker\u2019s pickups. A scoring attribute receives +1 point, while a penalty attribute receives -1 point. A higher total score for a picker denotes a higher quality of pickups.The table below summarizes the assessment of four pickers at your warehouse during one workday.Which picker scored highest in pickup quality on this day?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"A`",`n            `"B`",`n            `"C`",`n            `"D`"`n        ],`n        `"score`": `"B`"`n    }`n]"

# Remove the top formatting row (A1 with value 0, bold font, border) by deleting the entire row
$ws.Rows.Item(1).Delete()

# Set the new text content (now shifted into A1)
$ws.Range("A1").Value = $newText
$ws.Rows.Item(1).AutoFit()
